$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de
# A new row for 91ade9cd... is inserted before ac4dd610...,
# and a new row for cf1ee4c2... is inserted between ac4dd610... and ed40f145....
# -----------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Hyperlinks.Delete()

$wsOverview.Range("A2").Value = "91ade9cd-6393-4687-b754-a2a72f29a25d.md"
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"

$wsOverview.Range("A3").Value = "ac4dd610-f844-4d1c-9897-5391f84bd420.md"
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

$wsOverview.Range("A4").Value = "cf1ee4c2-fab2-4202-bc68-8ffe786fd0bc.md"
$wsOverview.Range("B4").Value = "Ready for handoff"
$wsOverview.Range("C4").Value = "Ready for handoff"

$wsOverview.Range("A5").Value = "ed40f145-c38d-42a2-8a32-92756b123b73.md"
$wsOverview.Range("B5").Value = "Ready for handoff"
$wsOverview.Range("C5").Value = "Ready for handoff"

$wsOverview.Range("A6").Value = ".localization-config"
$wsOverview.Range("B6").Value = "Not to be localized"
$wsOverview.Range("C6").Value = "Not to be localized"

$wsOverview.Range("A2:A6").Style = "HyperLink"

$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/590091b1b02240f3223bcfb53d28b23d3a15f259/e2e/91ade9cd-6393-4687-b754-a2a72f29a25d.md", [Type]::Missing, [Type]::Missing, "91ade9cd-6393-4687-b754-a2a72f29a25d.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/590091b1b02240f3223bcfb53d28b23d3a15f259/e2e/ac4dd610-f844-4d1c-9897-5391f84bd420.md", [Type]::Missing, [Type]::Missing, "ac4dd610-f844-4d1c-9897-5391f84bd420.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/590091b1b02240f3223bcfb53d28b23d3a15f259/e2e/cf1ee4c2-fab2-4202-bc68-8ffe786fd0bc.md", [Type]::Missing, [Type]::Missing, "cf1ee4c2-fab2-4202-bc68-8ffe786fd0bc.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/590091b1b02240f3223bcfb53d28b23d3a15f259/e2e/ed40f145-c38d-42a2-8a32-92756b123b73.md", [Type]::Missing, [Type]::Missing, "ed40f145-c38d-42a2-8a32-92756b123b73.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/590091b1b02240f3223bcfb53d28b23d3a15f259/.localization-config", [Type]::Missing, [Type]::Missing, ".localization-config")

# -----------------------------------------------------------------
# Sheet "zh-cn": Source File Name | Status | Latest Handoff File | Latest Handoff Datetime |
#                Latest Target File | Latest Handback File | Latest Handback DateTime |
#                Handoff Reason | Dependency From
# -----------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Hyperlinks.Delete()

$wsZh.Range("A2").Value = "91ade9cd-6393-4687-b754-a2a72f29a25d.md"
$wsZh.Range("B2").Value = "Ready for handoff"
$wsZh.Range("C2").Value = "91ade9cd-6393-4687-b754-a2a72f29a25d.3ba27fb5ec903485f7f4a11503412a1ddd36066f.zh-cn.xlf"
$wsZh.Range("D2").Value = "2016-03-08 01:51:41"
$wsZh.Range("G2").Value = "0001-01-01 00:00:00"
$wsZh.Range("H2").Value = "Include"

$wsZh.Range("A3").Value = "ac4dd610-f844-4d1c-9897-5391f84bd420.md"
$wsZh.Range("B3").Value = "Ready for handoff"
$wsZh.Range("C3").Value = "ac4dd610-f844-4d1c-9897-5391f84bd420.69d3935aed3c9c76cd34014921001c2b364c62e4.zh-cn.xlf"
$wsZh.Range("D3").Value = "2016-03-08 01:51:06"
$wsZh.Range("G3").Value = "0001-01-01 00:00:00"
$wsZh.Range("H3").Value = "Include"

$wsZh.Range("A4").Value = "cf1ee4c2-fab2-4202-bc68-8ffe786fd0bc.md"
$wsZh.Range("B4").Value = "Ready for handoff"
$wsZh.Range("C4").Value = "cf1ee4c2-fab2-4202-bc68-8ffe786fd0bc.f13c651d35e85c0337e8ad1958fd61cf3e0b08e5.zh-cn.xlf"
$wsZh.Range("D4").Value = "2016-03-08 01:51:41"
$wsZh.Range("G4").Value = "0001-01-01 00:00:00"
$wsZh.Range("H4").Value = "Include"

$wsZh.Range("A5").Value = "ed40f145-c38d-42a2-8a32-92756b123b73.md"
$wsZh.Range("B5").Value = "Ready for handoff"
$wsZh.Range("C5").Value = "ed40f145-c38d-42a2-8a32-92756b123b73.4c13a67f8404fcbac064efc971e47c3c7192f100.zh-cn.xlf"
$wsZh.Range("D5").Value = "2016-03-08 01:51:06"
$wsZh.Range("G5").Value = "0001-01-01 00:00:00"
$wsZh.Range("H5").Value = "Include"

$wsZh.Range("A6").Value = ".localization-config"
$wsZh.Range("B6").Value = "Not to be localized"
$wsZh.Range("D6").Value = "0001-01-01 00:00:00"
$wsZh.Range("G6").Value = "0001-01-01 00:00:00"
$wsZh.Range("H6").Value = "Ignored"

$wsZh.Range("A2:A6").Style = "HyperLink"
$wsZh.Range("C2").Style = "HyperLink"
$wsZh.Range("C3").Style = "HyperLink"
$wsZh.Range("C4").Style = "HyperLink"
$wsZh.Range("C5").Style = "HyperLink"
$wsZh.Range("D2:D6").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/590091b1b02240f3223bcfb53d28b23d3a15f259/e2e/91ade9cd-6393-4687-b754-a2a72f29a25d.md", [Type]::Missing, [Type]::Missing, "91ade9cd-6393-4687-b754-a2a72f29a25d.md")
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ced94659784d07d86d970e76db9ad38cd71a4fd7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/91ade9cd-6393-4687-b754-a2a72f29a25d.3ba27fb5ec903485f7f4a11503412a1ddd36066f.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "91ade9cd-6393-4687-b754-a2a72f29a25d.3ba27fb5ec903485f7f4a11503412a1ddd36066f.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/590091b1b02240f3223bcfb53d28b23d3a15f259/e2e/ac4dd610-f844-4d1c-9897-5391f84bd420.md", [Type]::Missing, [Type]::Missing, "ac4dd610-f844-4d1c-9897-5391f84bd420.md")
$wsZh.Hyperlinks.Add($wsZh.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ced94659784d07d86d970e76db9ad38cd71a4fd7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/ac4dd610-f844-4d1c-9897-5391f84bd420.69d3935aed3c9c76cd34014921001c2b364c62e4.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "ac4dd610-f844-4d1c-9897-5391f84bd420.69d3935aed3c9c76cd34014921001c2b364c62e4.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/590091b1b02240f3223bcfb53d28b23d3a15f259/e2e/cf1ee4c2-fab2-4202-bc68-8ffe786fd0bc.md", [Type]::Missing, [Type]::Missing, "cf1ee4c2-fab2-4202-bc68-8ffe786fd0bc.md")
$wsZh.Hyperlinks.Add($wsZh.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ced94659784d07d86d970e76db9ad38cd71a4fd7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/cf1ee4c2-fab2-4202-bc68-8ffe786fd0bc.f13c651d35e85c0337e8ad1958fd61cf3e0b08e5.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "cf1ee4c2-fab2-4202-bc68-8ffe786fd0bc.f13c651d35e85c0337e8ad1958fd61cf3e0b08e5.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/590091b1b02240f3223bcfb53d28b23d3a15f259/e2e/ed40f145-c38d-42a2-8a32-92756b123b73.md", [Type]::Missing, [Type]::Missing, "ed40f145-c38d-42a2-8a32-92756b123b73.md")
$wsZh.Hyperlinks.Add($wsZh.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ced94659784d07d86d970e76db9ad38cd71a4fd7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/ed40f145-c38d-42a2-8a32-92756b123b73.4c13a67f8404fcbac064efc971e47c3c7192f100.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "ed40f145-c38d-42a2-8a32-92756b123b73.4c13a67f8404fcbac064efc971e47c3c7192f100.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/590091b1b02240f3223bcfb53d28b23d3a15f259/.localization-config", [Type]::Missing, [Type]::Missing, ".localization-config")

# -----------------------------------------------------------------
# Sheet "de-de": same columns as zh-cn, but de-de target files/dates.
# -----------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Hyperlinks.Delete()

$wsDe.Range("A2").Value = "91ade9cd-6393-4687-b754-a2a72f29a25d.md"
$wsDe.Range("B2").Value = "Ready for handoff"
$wsDe.Range("C2").Value = "91ade9cd-6393-4687-b754-a2a72f29a25d.3ba27fb5ec903485f7f4a11503412a1ddd36066f.de-de.xlf"
$wsDe.Range("D2").Value = "2016-03-08 01:51:48"
$wsDe.Range("G2").Value = "0001-01-01 00:00:00"
$wsDe.Range("H2").Value = "Include"

$wsDe.Range("A3").Value = "ac4dd610-f844-4d1c-9897-5391f84bd420.md"
$wsDe.Range("B3").Value = "Ready for handoff"
$wsDe.Range("C3").Value = "ac4dd610-f844-4d1c-9897-5391f84bd420.69d3935aed3c9c76cd34014921001c2b364c62e4.de-de.xlf"
$wsDe.Range("D3").Value = "2016-03-08 01:51:25"
$wsDe.Range("G3").Value = "0001-01-01 00:00:00"
$wsDe.Range("H3").Value = "Include"

$wsDe.Range("A4").Value = "cf1ee4c2-fab2-4202-bc68-8ffe786fd0bc.md"
$wsDe.Range("B4").Value = "Ready for handoff"
$wsDe.Range("C4").Value = "cf1ee4c2-fab2-4202-bc68-8ffe786fd0bc.f13c651d35e85c0337e8ad1958fd61cf3e0b08e5.de-de.xlf"
$wsDe.Range("D4").Value = "2016-03-08 01:51:48"
$wsDe.Range("G4").Value = "0001-01-01 00:00:00"
$wsDe.Range("H4").Value = "Include"

$wsDe.Range("A5").Value = "ed40f145-c38d-42a2-8a32-92756b123b73.md"
$wsDe.Range("B5").Value = "Ready for handoff"
$wsDe.Range("C5").Value = "ed40f145-c38d-42a2-8a32-92756b123b73.4c13a67f8404fcbac064efc971e47c3c7192f100.de-de.xlf"
$wsDe.Range("D5").Value = "2016-03-08 01:51:25"
$wsDe.Range("G5").Value = "0001-01-01 00:00:00"
$wsDe.Range("H5").Value = "Include"

$wsDe.Range("A6").Value = ".localization-config"
$wsDe.Range("B6").Value = "Not to be localized"
$wsDe.Range("D6").Value = "0001-01-01 00:00:00"
$wsDe.Range("G6").Value = "0001-01-01 00:00:00"
$wsDe.Range("H6").Value = "Ignored"

$wsDe.Range("A2:A6").Style = "HyperLink"
$wsDe.Range("C2").Style = "HyperLink"
$wsDe.Range("C3").Style = "HyperLink"
$wsDe.Range("C4").Style = "HyperLink"
$wsDe.Range("C5").Style = "HyperLink"
$wsDe.Range("D2:D6").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/590091b1b02240f3223bcfb53d28b23d3a15f259/e2e/91ade9cd-6393-4687-b754-a2a72f29a25d.md", [Type]::Missing, [Type]::Missing, "91ade9cd-6393-4687-b754-a2a72f29a25d.md")
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5099fe94d844453caae6ca54ea2a9aa2da9b9166/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/91ade9cd-6393-4687-b754-a2a72f29a25d.3ba27fb5ec903485f7f4a11503412a1ddd36066f.de-de.xlf", [Type]::Missing, [Type]::Missing, "91ade9cd-6393-4687-b754-a2a72f29a25d.3ba27fb5ec903485f7f4a11503412a1ddd36066f.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/590091b1b02240f3223bcfb53d28b23d3a15f259/e2e/ac4dd610-f844-4d1c-9897-5391f84bd420.md", [Type]::Missing, [Type]::Missing, "ac4dd610-f844-4d1c-9897-5391f84bd420.md")
$wsDe.Hyperlinks.Add($wsDe.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5099fe94d844453caae6ca54ea2a9aa2da9b9166/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/ac4dd610-f844-4d1c-9897-5391f84bd420.69d3935aed3c9c76cd34014921001c2b364c62e4.de-de.xlf", [Type]::Missing, [Type]::Missing, "ac4dd610-f844-4d1c-9897-5391f84bd420.69d3935aed3c9c76cd34014921001c2b364c62e4.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/590091b1b02240f3223bcfb53d28b23d3a15f259/e2e/cf1ee4c2-fab2-4202-bc68-8ffe786fd0bc.md", [Type]::Missing, [Type]::Missing, "cf1ee4c2-fab2-4202-bc68-8ffe786fd0bc.md")
$wsDe.Hyperlinks.Add($wsDe.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5099fe94d844453caae6ca54ea2a9aa2da9b9166/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/cf1ee4c2-fab2-4202-bc68-8ffe786fd0bc.f13c651d35e85c0337e8ad1958fd61cf3e0b08e5.de-de.xlf", [Type]::Missing, [Type]::Missing, "cf1ee4c2-fab2-4202-bc68-8ffe786fd0bc.f13c651d35e85c0337e8ad1958fd61cf3e0b08e5.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/590091b1b02240f3223bcfb53d28b23d3a15f259/e2e/ed40f145-c38d-42a2-8a32-92756b123b73.md", [Type]::Missing, [Type]::Missing, "ed40f145-c38d-42a2-8a32-92756b123b73.md")
$wsDe.Hyperlinks.Add($wsDe.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5099fe94d844453caae6ca54ea2a9aa2da9b9166/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/ed40f145-c38d-42a2-8a32-92756b123b73.4c13a67f8404fcbac064efc971e47c3c7192f100.de-de.xlf", [Type]::Missing, [Type]::Missing, "ed40f145-c38d-42a2-8a32-92756b123b73.4c13a67f8404fcbac064efc971e47c3c7192f100.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/590091b1b02240f3223bcfb53d28b23d3a15f259/.localization-config", [Type]::Missing, [Type]::Missing, ".localization-config")

$wb.Worksheets.Item("Overview").Activate()
